$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the cell values can be
# updated, then restore protection once the edits are complete.
$ws.Unprotect()

# Update the "as of" date in the confidential disclosure text (cell A9)
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-27 for illustrative purposes only and are subject to change."

# Update weight / percent-change figures in rows 2-6
$ws.Range("D2").Value = 0.2536103295502961
$ws.Range("E2").Value = 0.007962150934687218

$ws.Range("D3").Value = 0.2532084991995142
$ws.Range("E3").Value = 0.01146972526006929

$ws.Range("D4").Value = 0.2433013806146668
$ws.Range("E4").Value = 0.01369334619093543

$ws.Range("D5").Value = 0.2498797906355229
$ws.Range("E5").Value = 0.003809523809523707

$ws.Range("E6").Value = 0.009207048687242203

# Restore sheet protection
$ws.Protect()
